$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")

# Make BOM the active sheet / active tab (moves activeTab from Sheet2 to BOM)
$ws.Activate()

# Update existing quantity/price values (C2, C3) - formulas in D2/D3 recalc automatically
$ws.Range("C2").Value = 67
$ws.Range("C3").Value = 120

# New rows: "Бачок" and "Шланг омывателя 5м"
$ws.Range("A4").Value = "Бачок"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 80

$ws.Range("A5").Value = "Шланг омывателя 5м"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 50

# Shared formula across D4:D5
$ws.Range("D4:D5").Formula = "=B4*C4"

# Totals row
$ws.Range("C8").Value = "Итого"
$ws.Range("D8").Formula = "=SUM(D2:D7)"

# Apply the "Calculation" built-in style to the amount column cells
$ws.Range("D2").Style = "Calculation"
$ws.Range("D3").Style = "Calculation"
$ws.Range("D4").Style = "Calculation"
$ws.Range("D5").Style = "Calculation"

# Update selection on BOM to D9, matching final cursor position
$ws.Range("D9").Select()
